$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = "Region"
$ws.Range("A2").Value = "us-east-1a"
$ws.Range("A3").Value = "us-east-1b"
$ws.Range("A4").Value = "us-west-1b"
